$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new row at position 6, pushing existing rows 6-12 down to 7-13
$ws.Rows.Item(6).Insert()

# Populate the new row with the "waittopageload1" / 2000 pair (T2 scenario row)
$ws.Range("A6").Value = "waittopageload1"
$ws.Range("B6").Value = 2000

# Match formatting of B6 to the existing similar row (B3), since Insert()
# otherwise carries over the format of the row above (row 5)
$ws.Range("B3").Copy()
$ws.Range("B6").PasteSpecial(-4122)

# Make "Edit Repayment Schedule" the active sheet/tab with A6:B6 selected
$ws.Activate()
$ws.Range("A6:B6").Select()
